$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: new review entry, formatted like the existing data rows ---
$ws.Range("A32:G32").Copy()
$ws.Range("A33:G33").PasteSpecial(-4122)

$ws.Range("A33").Value2 = "com.hamxa.shaynachim"
$ws.Range("B33").Value2 = "bitcoin"
$ws.Range("C33").Value2 = "shamirnaftali@gmail.com"
$ws.Range("D33").Value2 = "irisalmog47@gmail.com"
$ws.Range("E33").Value2 = "27/8/2019 19:59"
$ws.Range("F33").Value2 = "so good so far"
$ws.Range("G33").Value2 = "no"

$ws.Rows.Item(33).RowHeight = 13.8

# --- Row 34: new review entry with a real date/time value + wrapped review text ---
$ws.Range("A32:G32").Copy()
$ws.Range("A34:G34").PasteSpecial(-4122)

$ws.Range("A34").Value2 = "com.hamxa.shaynachim"
$ws.Range("B34").Value2 = "bitcoin"
$ws.Range("C34").Value2 = "kevinkors122@gmail.com"
$ws.Range("D34").Value2 = "sinuspai@gmail.com"

$ws.Range("E34").NumberFormat = "MM/DD/YYYY\ HH:MM:SS"
$ws.Range("E34").Value2 = 43746.1243055556

$ws.Range("F34").WrapText = $true
$ws.Range("F34").Value2 = "Great app, has a lot of info on crypto "

$ws.Range("G34").Value2 = "yes"

$ws.Rows.Item(34).RowHeight = 12.8

# --- Hyperlinks for the new emails in row 34 (row 33 stays plain text) ---
$ws.Hyperlinks.Add($ws.Range("C34"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D34"), "mailto:sinuspai@gmail.com", "", "", "sinuspai@gmail.com")

# --- Keep the view pointed at the new last row, like the source workbook ---
$ws.Range("E34").Select()
$excel.ActiveWindow.ScrollRow = 19
